$values = @(
    'Jumpsuit,Kaftan',
    'Kaftan,Tee',
    'Jumpsuit,Halter',
    'Turtleneck,Halter',
    'Trunks,Cutoffs',
    'Jumpsuit,Halter',
    'Blouse,Trunks',
    'Cutoffs,Trunks',
    'Jumpsuit,Blouse',
    'Jumpsuit,Trunks',
    'Halter,Blazer',
    'Jumpsuit,Blouse',
    'Jodhpurs,Trunks',
    'Jumpsuit,Blouse',
    'Dress,Jumpsuit',
    'Jumpsuit,Kaftan',
    'Jumpsuit,Dress',
    'Jumpsuit,Dress',
    'Jumpsuit,Caftan',
    'Jumpsuit,Kaftan',
    'Trunks,Blouse',
    'Jumpsuit,Blouse',
    'Trunks,Kaftan',
    'Kaftan,Blouse',
    'Jumpsuit,Blouse',
    'Jumpsuit,Kaftan',
    'Halter,Trunks',
    'Jumpsuit,Kaftan',
    'Jumpsuit,Halter',
    'Jumpsuit,Kaftan',
    'Jumpsuit,Blouse',
    'Jumpsuit,Kaftan',
    'Parka,Caftan',
    'Parka,Dress',
    'Jumpsuit,Halter',
    'Blazer,Halter',
    'Blazer,Halter',
    'Halter,Blazer',
    'Halter,Caftan',
    'Jumpsuit,Chinos',
    'Blouse,Jumpsuit',
    'Jodhpurs,Trunks',
    'Jumpsuit,Trunks',
    'Jumpsuit,Parka',
    'Jumpsuit,Blouse',
    'Jumpsuit,Blouse',
    'Jumpsuit,Blouse',
    'Trunks,Jumpsuit',
    'Jumpsuit,Jodhpurs',
    'Jumpsuit,Tee',
    'Kaftan,Jumpsuit',
    'Jumpsuit,Tee',
    'Jumpsuit,Trunks',
    'Jumpsuit,Parka',
    'Kaftan,Jumpsuit',
    'Jumpsuit,Parka',
    'Parka,Jumpsuit',
    'Jumpsuit,Dress',
    'Jumpsuit,Dress',
    'Halter,Caftan',
    'Jumpsuit,Kaftan',
    'Jumpsuit,Dress',
    'Jumpsuit,Kaftan',
    'Jumpsuit,Blouse',
    'Turtleneck,Jumpsuit',
    'Jumpsuit,Parka',
    'Jumpsuit,Kaftan',
    'Jumpsuit,Kaftan',
    'Jumpsuit,Dress',
    'Jumpsuit,Blouse',
    'Jumpsuit,Kaftan',
    'Jumpsuit,Halter',
    'Jumpsuit,Trunks',
    'Jumpsuit,Jodhpurs',
    'Jumpsuit,Halter',
    'Jumpsuit,Blouse',
    'Jumpsuit,Trunks',
    'Jumpsuit,Dress',
    'Trunks,Blouse',
    'Jodhpurs,Jumpsuit',
    'Halter,Blazer',
    'Kaftan,Tee',
    'Jumpsuit,Halter',
    'Halter,Jumpsuit',
    'Jumpsuit,Halter',
    'Jumpsuit,Dress',
    'Jumpsuit,Tee',
    'Jumpsuit,Kaftan',
    'Jumpsuit,Dress',
    'Trunks,Blazer',
    'Jumpsuit,Halter',
    'Jumpsuit,Blazer',
    'Jumpsuit,Dress',
    'Jumpsuit,Kaftan',
    'Caftan,Trunks',
    'Jumpsuit,Halter'
)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $values[$i]
}
